# Glossary.xlsx -- References sheet: add new glossary/reference rows (4-18),
# their hyperlinks, and resize columns to fit the new (longer) URL content.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Set cell text values in the order the original author entered them --
#    this drives shared-string allocation order so it matches the target file.
$ws.Range("A4").Value = "Adnodd codio cymraeg "
$ws.Range("B4").Value = "Adnodd data/maths uwch yn gymraeg"
$ws.Range("A5").Value = "https://adnoddau.porth.ac.uk/webapps/portal/execute/tabs/tabAction?tab_tab_group_id=_90_1"
$ws.Range("A6").Value = "https://projects.raspberrypi.org/cy-GB/codeclub"
$ws.Range("A7").Value = "https://www.technocamps.com/cy/resources/python"
$ws.Range("A8").Value = "https://github.com/meigwilym/haciaith13"
$ws.Range("A10").Value = "https://sgiliauymchwilcyfrifiadurol.github.io/"
$ws.Range("B5").Value = "https://llyfrgell.porth.ac.uk/View.aspx?id=1716~4p~QbzBunJu"
$ws.Range("A9").Value = "https://github.com/porthtechnolegauiaith"
$ws.Range("A11").Value = "http://techiaith.cymru/yr-adnoddau/llawlyfr-technolegau-iaith/"
$ws.Range("A12").Value = "https://www.meddal.com/meddal/?page_id=1111"
$ws.Range("A13").Value = "https://www.technocamps.com/cy/resources/arduino"
$ws.Range("A14").Value = "https://apps.apple.com/us/app/botio/id1296278646?ls=1"
$ws.Range("A15").Value = "https://www.technocamps.com/cy/resources/artificial-intelligence"
$ws.Range("A16").Value = "https://www.technocamps.com/cy/resources"
$ws.Range("A17").Value = "http://resources.hwb.wales.gov.uk/VTC/ngfl/computing/164/index.html"
$ws.Range("B6").Value = "https://www.mathemateg.com/"
$ws.Range("A18").Value = "https://llyfrgell.porth.ac.uk/Default.aspx?search=python&pagesize=20&page=2&fp=0"
$ws.Range("B8").Value = "https://llyfrgell.porth.ac.uk/View.aspx?id=5729~4x~8AS68GtY"
$ws.Range("B9").Value = "https://adnoddau.cbac.co.uk/Pages/ResourceSingle.aspx?rIid=2660"

# 2) B7 is an empty cell that only carries the Hyperlink style (no text/link).
$ws.Range("B7").Style = "Hyperlink"

# 3) Wire up the hyperlink relationships (row by row, top to bottom).
$ws.Hyperlinks.Add($ws.Range("A5"), "https://adnoddau.porth.ac.uk/webapps/portal/execute/tabs/tabAction?tab_tab_group_id=_90_1") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "https://projects.raspberrypi.org/cy-GB/codeclub") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), "https://www.technocamps.com/cy/resources/python") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A8"), "https://github.com/meigwilym/haciaith13") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A9"), "https://github.com/porthtechnolegauiaith") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A10"), "https://sgiliauymchwilcyfrifiadurol.github.io/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "https://llyfrgell.porth.ac.uk/View.aspx?id=1716~4p~QbzBunJu") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A11"), "http://techiaith.cymru/yr-adnoddau/llawlyfr-technolegau-iaith/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A12"), "https://www.meddal.com/meddal/?page_id=1111") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A13"), "https://www.technocamps.com/cy/resources/arduino") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A14"), "https://apps.apple.com/us/app/botio/id1296278646?ls=1") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A15"), "https://www.technocamps.com/cy/resources/artificial-intelligence") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A16"), "https://www.technocamps.com/cy/resources") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A17"), "http://resources.hwb.wales.gov.uk/VTC/ngfl/computing/164/index.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B6"), "https://www.mathemateg.com/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A18"), "https://llyfrgell.porth.ac.uk/Default.aspx?search=python&pagesize=20&page=2&fp=0") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B8"), "https://llyfrgell.porth.ac.uk/View.aspx?id=5729~4x~8AS68GtY") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B9"), "https://adnoddau.cbac.co.uk/Pages/ResourceSingle.aspx?rIid=2660") | Out-Null

# 4) Re-apply the Hyperlink style so every linked cell reuses cell style index 1
#    (the same style already used by B1/B2), instead of a freshly minted one.
$ws.Range("A5").Style = "Hyperlink"
$ws.Range("A6").Style = "Hyperlink"
$ws.Range("A7").Style = "Hyperlink"
$ws.Range("A8").Style = "Hyperlink"
$ws.Range("A9").Style = "Hyperlink"
$ws.Range("A10").Style = "Hyperlink"
$ws.Range("B5").Style = "Hyperlink"
$ws.Range("A11").Style = "Hyperlink"
$ws.Range("A12").Style = "Hyperlink"
$ws.Range("A13").Style = "Hyperlink"
$ws.Range("A14").Style = "Hyperlink"
$ws.Range("A15").Style = "Hyperlink"
$ws.Range("A16").Style = "Hyperlink"
$ws.Range("A17").Style = "Hyperlink"
$ws.Range("B6").Style = "Hyperlink"
$ws.Range("A18").Style = "Hyperlink"
$ws.Range("B8").Style = "Hyperlink"
$ws.Range("B9").Style = "Hyperlink"

# 5) Resize the two columns so the much longer URLs fit (author re-ran AutoFit).
$ws.Columns.Item(1).ColumnWidth = 88.66666666666667
$ws.Columns.Item(2).ColumnWidth = 178.0

# 6) Restore the selection to where the author left off editing.
$ws.Range("A15").Select() | Out-Null
